$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cyclones")

# Insert a new blank row above the "TEAM TOTALS" row (row 13), pushing
# TEAM TOTALS down to row 14 and growing the sheet's used range by one row.
$ws.Rows.Item(13).Insert()

# Give the (now) row 14 "TEAM TOTALS" row a thin box border around every
# cell, matching the look of the bold/bordered header row above the data.
$totalsRow = $ws.Range("A14:P14")
$totalsRow.Borders.LineStyle = 1
